$wb = $excel.ActiveWorkbook

# Add the new "Wyvern" sheet after the last existing sheet (Beetle)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Wyvern"

# Header row
$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "D20"
$ws.Range("C1").Value = "Armor"
$ws.Range("D1").Value = "HP"

# Hit location data rows
$data = @(
    @("Tail",       1,        6, 6),
    @("Right Leg",  "02-04",  6, 6),
    @("Left Leg",   "05-07",  6, 6),
    @("Abdomen",    "08-08",  6, 6),
    @("Chest",      "09-11",  6, 7),
    @("Right Wing", "13-14",  6, 5),
    @("Left Wing",  "15-16",  6, 5),
    @("Head",       "17-20",  6, 6)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the saved view state: Dragonsnail1 tab active, with C25 selected
$ws2 = $wb.Worksheets.Item("Dragonsnail1")
$ws2.Activate() | Out-Null
$ws2.Range("C25").Select() | Out-Null

$wb.Save()
